# Actualización SmartScore desde Streamlit (Miranda)
#
# 1. Copy the existing row 5 (Miranda's first submission) down to a new
#    row 6, preserving cell types/formatting exactly (Copy/PasteSpecial
#    keeps text cells as text instead of Excel re-interpreting numeric
#    looking strings as numbers).
# 2. Update the timestamp in the new row (D6) to the later submission time.
# 3. Fix up row 5 so the nine SmartScore columns (G, J, M, P, S, V, Y, AB,
#    AE), which were stored as text in the source file, become real
#    numeric values - matching how every other row in the sheet stores
#    SmartScore.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: duplicate row 5 into row 6 -----------------------------------
$ws.Range("A5:AF5").Copy($ws.Range("A6:AF6"))

# --- Step 2: new row gets its own submission timestamp ---------------------
$ws.Cells.Item(6, 4).Value = "2025-10-28 05:56:46"

# --- Step 3: convert row 5's SmartScore text cells into numbers -----------
$scoreCols = @(7, 10, 13, 16, 19, 22, 25, 28, 31)   # G, J, M, P, S, V, Y, AB, AE
$scoreVals = @(0.575, 0.51, 0.509, 0.65, 0.587, 0.552, 0.664, 0.589, 0.576)

for ($i = 0; $i -lt $scoreCols.Length; $i++) {
    $ws.Cells.Item(5, $scoreCols[$i]).Value = $scoreVals[$i]
}
